$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.820.38'
$ws.Range("E2").Value = '  +0.86%  '

$ws.Range("D3").Value = '1.810.03'
$ws.Range("E3").Value = '  +0.49%  '

$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.19'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.05%  '

$ws.Range("E6").Value = '  +0.36%  '

$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.15'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.42%  '

$ws.Range("E9").Value = '  +4.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0679'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0999'
$ws.Range("D11").ClearFormats()

$ws.Range("D12").Value = '2.072.26'
$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("D13").Value = '1.810.40'
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.05'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.29%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.655'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.17%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.63'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +5.29%  '

$ws.Range("D17").Value = '34.810.06'
$ws.Range("E17").Value = '  +0.98%  '

$ws.Range("E18").Value = '  +2.42%  '

$ws.Range("D19").Value = '0.0₃0781'
$ws.Range("E19").Value = '  +1.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.70%  '

$ws.Range("E21").Value = '  +5.26%  '

$ws.Range("E22").Value = '  +9.41%  '

$ws.Range("E23").Value = '  +0.37%  '

$ws.Range("E24").Value = '  +4.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.28'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.70'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.29'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.119'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.99%  '

$ws.Range("E29").Value = '  +29.10%  '

$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("D31").Value = '3.344.10'
$ws.Range("E31").Value = '  +37.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0541'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.85'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.56%  '

$ws.Range("E34").Value = '  +1.99%  '

$ws.Range("E35").Value = '  -1.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '92.93'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.39%  '

$ws.Range("E37").Value = '  +5.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.672'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.88%  '

$ws.Range("D39").Value = '1.302.20'
$ws.Range("E39").Value = '  -1.60%  '

$ws.Range("E40").Value = '  +4.10%  '

$ws.Range("E41").Value = '  +1.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.90'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.977'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.55%  '

$ws.Range("E44").Value = '  -0.63%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.74'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.15'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0513'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.75%  '

$ws.Range("D49").Value = '1.986.93'
$ws.Range("E49").Value = '  +1.16%  '

$ws.Range("E50").Value = '  +0.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.57'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.08%  '
